$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2,3,4 data gets cyclically rotated: new row2 = old row4, new row3 = old row2, new row4 = old row3.
# Only columns D, J, K, L, M, P change; everything else stays the same.

# Row 2 (was old row 4)
$ws.Range("D2").Value = 44839
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15600
$ws.Range("P2").Value = 1040

# Row 3 (was old row 2)
$ws.Range("D3").Value = 44749
$ws.Range("J3").Value = 90
$ws.Range("K3").Value = 17000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 17556
$ws.Range("P3").Value = 1170

# Row 4 (was old row 3)
$ws.Range("D4").Value = 44750
$ws.Range("J4").Value = 140
$ws.Range("K4").Value = 19000
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 19571
$ws.Range("P4").Value = 1305
